# PDF export functionality added + some refactoring.
#
# This script:
#  1) Bumps Сидорова's Математика (Maths) grade on the raw data sheet from
#     3 to 4, and propagates the same value to its mirror on the per-subject
#     "Математика" sheet, plus the recalculated per-student / overall
#     average rows (stored as static values in this workbook, so each is
#     set explicitly).
#  2) Repositions/resizes the four embedded charts (one per chart sheet) so
#     they anchor from row 11 / col A down to around row 26, instead of the
#     original row 2 / col D band.

$wb = $excel.ActiveWorkbook

$EMU_PER_POINT = 12700

# Anchors a ChartObject's top-left/bottom-right to precise (0-indexed) OOXML
# col/colOff/row/rowOff coordinates, the same way <xdr:from>/<xdr:to> encode
# a twoCellAnchor, by driving Left/Top/Width/Height off real cell geometry.
function Set-ChartAnchor {
    param(
        $Chart,
        $Sheet,
        [int]$FromCol, [double]$FromColOff, [int]$FromRow, [double]$FromRowOff,
        [int]$ToCol,   [double]$ToColOff,   [int]$ToRow,   [double]$ToRowOff
    )

    # xdr uses 0-indexed col/row; Cells.Item is 1-indexed, so +1.
    $fromCell = $Sheet.Cells.Item($FromRow + 1, $FromCol + 1)
    $toCell   = $Sheet.Cells.Item($ToRow + 1, $ToCol + 1)

    $left = $fromCell.Left + ($FromColOff / $EMU_PER_POINT)
    $top  = $fromCell.Top  + ($FromRowOff / $EMU_PER_POINT)

    $right  = $toCell.Left + ($ToColOff / $EMU_PER_POINT)
    $bottom = $toCell.Top  + ($ToRowOff / $EMU_PER_POINT)

    $Chart.Left   = $left
    $Chart.Top    = $top
    $Chart.Width  = $right - $left
    $Chart.Height = $bottom - $top
}

# --- 1) Grade data -----------------------------------------------------

$wsGrades = $wb.Worksheets.Item("Успеваемость")
$wsGrades.Range("H8").Value = 4          # Сидорова, Математика: 3 -> 4
$wsGrades.Range("B16").Value = 4         # Сидорова average: 3.67 -> 4
$wsGrades.Range("B17").Value = 3.94      # Overall average: 3.89 -> 3.94

$wsAvg = $wb.Worksheets.Item("Диаграмма средней успеваемости")
$wsAvg.Range("B8").Value = 4             # mirrors Успеваемость!B16
$wsAvg.Range("B9").Value = 3.94          # mirrors Успеваемость!B17

$wsMath = $wb.Worksheets.Item("Математика  диаграмма успеваемо")
$wsMath.Range("B8").Value = 4            # mirrors Успеваемость!H8

# --- 2) Chart repositioning ---------------------------------------------

$chartSheets = @(
    "Диаграмма средней успеваемости",
    "Информатика  диаграмма успеваем",
    "История  диаграмма успеваемости",
    "Математика  диаграмма успеваемо"
)

$toTargets = @{
    "Диаграмма средней успеваемости"  = @(7, 171450, 25, 0)
    "Информатика  диаграмма успеваем" = @(8, 314325, 25, 0)
    "История  диаграмма успеваемости" = @(8, 314325, 25, 0)
    "Математика  диаграмма успеваемо" = @(8, 314325, 25, 0)
}

foreach ($sheetName in $chartSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $co = $ws.ChartObjects().Item(1)
    $to = $toTargets[$sheetName]

    # NOTE: this runtime's PowerShell only binds positional args reliably,
    # so call Set-ChartAnchor positionally (Chart, Sheet, FromCol, FromColOff,
    # FromRow, FromRowOff, ToCol, ToColOff, ToRow, ToRowOff).
    Set-ChartAnchor $co $ws 0 0 10 0 $to[0] $to[1] $to[2] $to[3]
}
